$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "68.059.88"
$cell.ClearFormats()
$ws.Cells.Item(2, 5).Value = "  -6.61%  "

# Row 3: Ethereum -> Ethereum
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.709.68"
$cell.ClearFormats()
$ws.Cells.Item(3, 5).Value = "  -6.11%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

# Row 5: BNB -> BNB
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "582.54"
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = "  -4.92%  "

# Row 6: Solana -> Solana
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "177.86"
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = "  +5.79%  "

# Row 7: LidoStakedEther -> LidoStakedEther
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.706.36"
$cell.ClearFormats()
$ws.Cells.Item(7, 5).Value = "  -5.91%  "

# Row 8: XRP -> XRP
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.635"
$cell.ClearFormats()
$ws.Cells.Item(8, 5).Value = "  -6.11%  "

# Row 9: USDC -> USDC
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.ClearFormats()
$ws.Cells.Item(9, 5).Value = "  -0.24%  "

# Row 10: Cardano -> Cardano
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.718"
$cell.ClearFormats()
$ws.Cells.Item(10, 5).Value = "  -4.16%  "

# Row 11: Dogecoin -> Dogecoin
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.166"
$cell.ClearFormats()
$ws.Cells.Item(11, 5).Value = "  -7.86%  "

# Row 12: Avalanche -> Avalanche
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "52.79"
$cell.ClearFormats()
$ws.Cells.Item(12, 5).Value = "  -5.69%  "

# Row 13: ShibaInu -> ShibaInu
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0000302"
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = "  -8.57%  "

# Row 14: Polkadot -> Polkadot
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.67"
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = "  -3.50%  "

# Row 15: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.302.34"
$cell.ClearFormats()
$ws.Cells.Item(15, 5).Value = "  -6.20%  "

# Row 16: WrappedEther -> WrappedEther
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.745.82"
$cell.ClearFormats()
$ws.Cells.Item(16, 5).Value = "  -5.52%  "

# Row 17: TRON -> Chainlink
$ws.Cells.Item(17, 2).Value = "Chainlink"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "19.45"
$cell.ClearFormats()
$ws.Cells.Item(17, 5).Value = "  -4.55%  "

# Row 18: Chainlink -> TRON
$ws.Cells.Item(18, 2).Value = "TRON"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.127"
$cell.ClearFormats()
$ws.Cells.Item(18, 5).Value = "  -2.99%  "

# Row 19: Uniswap -> Polygon
$ws.Cells.Item(19, 2).Value = "Polygon"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.14"
$cell.ClearFormats()
$ws.Cells.Item(19, 5).Value = "  -7.84%  "

# Row 20: Polygon -> Uniswap
$ws.Cells.Item(20, 2).Value = "Uniswap"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "13.02"
$cell.ClearFormats()
$ws.Cells.Item(20, 5).Value = "  -6.67%  "

# Row 21: WrappedBTC -> WrappedBTC
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "68.008.06"
$cell.ClearFormats()
$ws.Cells.Item(21, 5).Value = "  -6.57%  "

# Row 22: BitcoinCash -> BitcoinCash
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "410.15"
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = "  -6.02%  "

# Row 23: PancakeSwap -> PancakeSwap
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.61"
$cell.ClearFormats()
$ws.Cells.Item(23, 5).Value = "  -5.51%  "

# Row 24: Litecoin -> Litecoin
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "88.71"
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = "  -6.78%  "

# Row 25: ImmutableX -> ImmutableX
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.10"
$cell.ClearFormats()
$ws.Cells.Item(25, 5).Value = "  -7.49%  "

# Row 26: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.89"
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = "  -8.37%  "

# Row 27: RenderToken -> RenderToken
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.76"
$cell.ClearFormats()
$ws.Cells.Item(27, 5).Value = "  -2.34%  "

# Row 28: Toncoin -> Toncoin
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.86"
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = "  -5.31%  "

# Row 29: LEO -> LEO
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.97"
$cell.ClearFormats()
$ws.Cells.Item(29, 5).Value = "  +0.30%  "

# Row 30: Filecoin -> Filecoin
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.57"
$cell.ClearFormats()
$ws.Cells.Item(30, 5).Value = "  -8.04%  "

# Row 31: NEARProtocol -> NEARProtocol
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.14"
$cell.ClearFormats()
$ws.Cells.Item(31, 5).Value = "  +2.22%  "

# Row 32: EthereumClassic -> EthereumClassic
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "33.07"
$cell.ClearFormats()
$ws.Cells.Item(32, 5).Value = "  -7.76%  "

# Row 33: Cosmos -> Cosmos
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.78"
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = "  -6.00%  "

# Row 34: InjectiveProtocol -> InjectiveProtocol
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "44.66"
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -6.18%  "

# Row 35: Hedera -> Hedera
$ws.Cells.Item(35, 5).Value = "  -8.22%  "

# Row 36: OKB -> OKB
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "66.15"
$cell.ClearFormats()
$ws.Cells.Item(36, 5).Value = "  -5.30%  "

# Row 37: PEPE -> PEPE
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0926"
$cell.ClearFormats()
$ws.Cells.Item(37, 5).Value = "  -9.26%  "

# Row 38: Bittensor -> Bittensor
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "594.28"
$cell.ClearFormats()
$ws.Cells.Item(38, 5).Value = "  -6.50%  "

# Row 39: TheGraph -> TheGraph
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.404"
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = "  -5.67%  "

# Row 40: Dai -> dogwifhat
$ws.Cells.Item(40, 2).Value = "dogwifhat"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.35"
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = "  +16.24%  "

# Row 41: FirstDigitalUSD -> Dai
$ws.Cells.Item(41, 2).Value = "Dai"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(41, 5).Value = "  +0.20%  "

# Row 42: dogwifhat -> FirstDigitalUSD
$ws.Cells.Item(42, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.ClearFormats()
$ws.Cells.Item(42, 5).Value = "  -0.22%  "

# Row 43: Kaspa -> Kaspa
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.137"
$cell.ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -6.20%  "

# Row 44: ThetaToken -> ThetaToken
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.10"
$cell.ClearFormats()
$ws.Cells.Item(44, 5).Value = "  -10.07%  "

# Row 45: VeChain -> VeChain
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0442"
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = "  -8.05%  "

# Row 46: Fetch.AI -> THORChain
$ws.Cells.Item(46, 2).Value = "THORChain"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.52"
$cell.ClearFormats()
$ws.Cells.Item(46, 5).Value = "  -11.26%  "

# Row 47: THORChain -> Fetch.AI
$ws.Cells.Item(47, 2).Value = "Fetch.AI"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.59"
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = "  -0.10%  "

# Row 48: Stellar -> Stellar
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.135"
$cell.ClearFormats()
$ws.Cells.Item(48, 5).Value = "  -8.39%  "

# Row 49: Maker -> Maker
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.753.23"
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = "  -3.66%  "

# Row 50: ApeXProtocol -> ApeXProtocol
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.16"
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = "  -7.03%  "

# Row 51: WEMIXToken -> WEMIXToken
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.68"
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = "  -14.15%  "
